# "Mejoras reporte de permanencia"
# - Rename the existing "Data" sheet to "Resumen".
# - Add a new "Detalle" sheet after it, with its own header row + autofilter.
# - Both sheets keep their own hidden _FilterDatabase defined name.
# - The new sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Resumen"

# --- Add the new "Detalle" sheet right after "Resumen" ----------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Detalle"

# --- Header row for Detalle (set in this order so the shared-string
#     table gets the same index order as the authored workbook) -------
$ws2.Range("B1").Value = "Proceso Id"
$ws2.Range("C1").Value = "Cantidad de veces que el documento ingresa a la unidad"
$ws2.Range("D1").Value = "Total tiempo de permanencia"
$ws2.Range("A1").Value = "Unidad origen"
$ws2.Range("A1:D1").Font.Bold = $true

# --- Column widths for Detalle -----------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 54.333333333333333
$ws2.Columns.Item(2).ColumnWidth = 24.0
$ws2.Columns.Item(3).ColumnWidth = 53.333333333333333
$ws2.Columns.Item(4).ColumnWidth = 29.0
$ws2.Columns.Item(5).ColumnWidth = 26.666666666666668

# --- Autofilter on the new sheet (matches the width of the old one) ---
$ws2.Range("A1:E1").AutoFilter()

# --- Rebuild the _FilterDatabase defined names so Detalle's entry
#     comes first, followed by Resumen's (both hidden) -----------------
$oldFilterName = $wb.Names.Item(1)
$oldFilterName.Delete()

$null = $ws2.Names.Add("_xlnm._FilterDatabase", "=Detalle!`$A`$1:`$E`$1")
$wb.Names.Item(1).Visible = $false

$null = $ws1.Names.Add("_xlnm._FilterDatabase", "=Resumen!`$A`$1:`$E`$1")
$wb.Names.Item(2).Visible = $false

# --- Make Detalle the active/selected sheet, with A2 selected ---------
$ws2.Activate()
$ws2.Range("A2").Select()
